$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for columns I (I0) and J (IF), rows 2-45 contain new data.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the existing headers (e.g. H1) by copying
# its formatting onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J, rows 2-45
$data = @(
    @(2, 8, 8),
    @(3, 9, 9),
    @(4, 6, 7),
    @(5, 9, 9),
    @(6, 8, 8),
    @(7, 8, 8),
    @(8, 8, 9),
    @(9, 8, 8),
    @(10, 8, 8),
    @(11, 9, 9),
    @(12, 6, 7),
    @(13, 8, 8),
    @(14, 9, 9),
    @(15, 8, 8),
    @(16, 9, 9),
    @(17, 8, 8),
    @(18, 8, 8),
    @(19, 11, 11),
    @(20, 9, 9),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 7, 7),
    @(24, 7, 7),
    @(25, 7, 7),
    @(26, 9, 9),
    @(27, 8, 8),
    @(28, 7, 7),
    @(29, 9, 9),
    @(30, 5, 5),
    @(31, 8, 8),
    @(32, 7, 7),
    @(33, 7, 7),
    @(34, 6, 6),
    @(35, 6, 6),
    @(36, 8, 8),
    @(37, 9, 9),
    @(38, 7, 7),
    @(39, 6, 6),
    @(40, 7, 7),
    @(41, 6, 6),
    @(42, 8, 8),
    @(43, 7, 7),
    @(44, 5, 5),
    @(45, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $if = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}
